# ------------------------------------------------------------------
# PlayerPerformance_3391.xlsx update
#   1. Insert a new "Player Info" sheet at the front with the
#      player's basic bio info.
#   2. On the existing "ODI Batting" sheet: rename MATCH_CARD_LINK ->
#      MATCH_CODE and replace the full scorecard URL with just the
#      numeric match code. Also drop the stray empty inline-string
#      placeholder cells that used to live in column B on a handful
#      of rows.
#   3. Same MATCH_CARD_LINK -> MATCH_CODE treatment on "ODI Bowling"
#      (column B there).
#   4. Append a brand-new "ODI Batting Extra" sheet with additional
#      per-innings batting detail.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. "Player Info" sheet (inserted as the first sheet)
# ------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $piHeaders.Count; $c++) {
    $cell = $playerInfo.Cells.Item(1, $c)
    $cell.Value = $piHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$piA2 = $playerInfo.Cells.Item(2, 1)
$piA2.NumberFormat = "@"
$piA2.Value = "3391"
$playerInfo.Cells.Item(2, 2).Value = "Kevin Joseph O''Brien"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Medium Fast"

# ------------------------------------------------------------------
# 2. "ODI Batting" - MATCH_CARD_LINK -> MATCH_CODE
# ------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Cells.Item(1, 4).Value = "MATCH_CODE"

$battingLastRow = $batting.UsedRange.Rows.Count
for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $batting.Cells.Item($r, 4)
    $val = [string]$cell.Value()
    if ($val -match "MatchCode=(\d+)") {
        $cell.NumberFormat = "@"
        $cell.Value = $matches[1]
    }

    # rows that used to carry a stray empty placeholder cell in col B
    $bCell = $batting.Cells.Item($r, 2)
    $bVal = $bCell.Value()
    if ($bVal -ne $null -and [string]$bVal -eq "") {
        $bCell.ClearContents()
    }
}

# ------------------------------------------------------------------
# 3. "ODI Bowling" - MATCH_CARD_LINK -> MATCH_CODE
# ------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Cells.Item(1, 2).Value = "MATCH_CODE"

$bowlingLastRow = $bowling.UsedRange.Rows.Count
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $bowling.Cells.Item($r, 2)
    $val = [string]$cell.Value()
    if ($val -match "MatchCode=(\d+)") {
        $cell.NumberFormat = "@"
        $cell.Value = $matches[1]
    }
}

# ------------------------------------------------------------------
# 4. "ODI Batting Extra" sheet (appended as the last sheet)
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

$exHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $exHeaders.Count; $c++) {
    $cell = $extra.Cells.Item(1, $c)
    $cell.Value = $exHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$exData = @(
  ,@(@{T='s';V='4284'}, @{T='n';V='5'}, @{T='s';V='1'},  @{T='s';V='0'}, @{T='s';V='2.02%'},  @{T='s';V='NO'})
  ,@(@{T='s';V='4285'}, @{T='n';V='5'}, @{T='s';V='7'},  @{T='s';V='1'}, @{T='s';V='36.76%'}, @{T='s';V='NO'})
  ,@(@{T='s';V='4291'}, @{T='n';V='5'}, @{T='s';V='3'},  @{T='s';V='3'}, @{T='s';V='19.27%'}, @{T='s';V='NO'})
  ,@(@{T='s';V='4295'}, $null,          $null,           $null,         $null,               @{T='s';V='NO'})
  ,@(@{T='s';V='4299'}, @{T='n';V='5'}, @{T='s';V='4'},  @{T='s';V='0'}, @{T='s';V='15.24%'}, @{T='s';V='NO'})
  ,@(@{T='s';V='4301'}, @{T='n';V='5'}, @{T='s';V='2'},  @{T='s';V='0'}, @{T='s';V='6.70%'},  @{T='s';V='NO'})
  ,@(@{T='s';V='4343'}, @{T='n';V='5'}, @{T='s';V='1'},  @{T='s';V='0'}, @{T='s';V='1.94%'},  @{T='s';V='NO'})
  ,@(@{T='s';V='4347'}, @{T='n';V='5'}, @{T='s';V='1'},  @{T='s';V='0'}, @{T='s';V='2.48%'},  @{T='s';V='NO'})
  ,@(@{T='s';V='4352'}, $null,          $null,           $null,         $null,               @{T='s';V='NO'})
  ,@(@{T='s';V='4391'}, @{T='n';V='5'}, @{T='s';V='0'},  @{T='s';V='0'}, @{T='s';V='2.22%'},  @{T='s';V='NO'})
  ,@(@{T='s';V='4394'}, @{T='n';V='5'}, @{T='s';V='4'},  @{T='s';V='1'}, @{T='s';V='13.08%'}, @{T='s';V='NO'})
  ,@(@{T='s';V='4397'}, @{T='n';V='5'}, @{T='s';V='4'},  @{T='s';V='0'}, @{T='s';V='10.34%'}, @{T='s';V='NO'})
  ,@(@{T='s';V='4426'}, $null,          $null,           $null,         $null,               @{T='s';V='NO'})
  ,@(@{T='s';V='4427'}, @{T='n';V='5'}, @{T='s';V='0'},  @{T='s';V='0'}, @{T='s';V='1.42%'},  @{T='s';V='NO'})
  ,@(@{T='s';V='4428'}, @{T='n';V='5'}, @{T='s';V='1'},  @{T='s';V='1'}, @{T='s';V='6.38%'},  @{T='s';V='NO'})
  ,@(@{T='s';V='4439'}, $null,          $null,           $null,         $null,               @{T='s';V='NO'})
  ,@(@{T='s';V='4442'}, @{T='n';V='2'}, @{T='s';V='0'},  @{T='s';V='0'}, @{T='s';V='0.88%'},  @{T='s';V='NO'})
  ,@(@{T='s';V='4444'}, @{T='n';V='2'}, @{T='s';V='0'},  @{T='s';V='0'}, @{T='s';V='0.37%'},  @{T='s';V='NO'})
  ,@(@{T='s';V='4446'}, @{T='n';V='2'}, @{T='s';V='0'},  @{T='s';V='0'}, @{T='s';V='0.39%'},  @{T='s';V='NO'})
  ,@(@{T='s';V='4468'}, @{T='n';V='1'}, @{T='s';V='0'},  @{T='s';V='0'}, $null,               @{T='s';V='NO'})
)

$rowNum = 2
foreach ($row in $exData) {
    for ($c = 1; $c -le 6; $c++) {
        $item = $row[$c - 1]
        $cell = $extra.Cells.Item($rowNum, $c)
        if ($item -eq $null) {
            $cell.NumberFormat = "@"
        } elseif ($item.T -eq "n") {
            $cell.Value = [double]$item.V
        } else {
            $cell.NumberFormat = "@"
            $cell.Value = $item.V
        }
    }
    $rowNum = $rowNum + 1
}
